{"js": "// Replace the date and the 25 division-problem answers in the table.\n// Each \"before\" text is unique within the document, so a simple\n// search-and-replace keyed on exact text is unambiguous and safe\n// regardless of replacement order (no \"after\" value collides with any\n// other \"before\" value).\nconst replacements = [\n  [\"2025-05-03 Saturday\", \"2025-05-04 Sunday\"],\n  [\"87\u00f72=43, 1\", \"51\u00f72=25, 1\"],\n  [\"33\u00f79=3, 6\", \"87\u00f79=9, 6\"],\n  [\"51\u00f75=10, 1\", \"16\u00f79=1, 7\"],\n  [\"47\u00f72=23, 1\", \"36\u00f78=4, 4\"],\n  [\"97\u00f72=48, 1\", \"25\u00f76=4, 1\"],\n  [\"71\u00f72=35, 1\", \"10\u00f78=1, 2\"],\n  [\"71\u00f75=14, 1\", \"18\u00f75=3, 3\"],\n  [\"55\u00f76=9, 1\", \"96\u00f75=19, 1\"],\n  [\"99\u00f75=19, 4\", \"38\u00f77=5, 3\"],\n  [\"40\u00f78=5, 0\", \"13\u00f79=1, 4\"],\n  [\"58\u00f72=29, 0\", \"23\u00f75=4, 3\"],\n  [\"13\u00f78=1, 5\", \"82\u00f73=27, 1\"],\n  [\"53\u00f76=8, 5\", \"17\u00f75=3, 2\"],\n  [\"35\u00f72=17, 1\", \"31\u00f73=10, 1\"],\n  [\"79\u00f77=11, 2\", \"52\u00f75=10, 2\"],\n  [\"36\u00f75=7, 1\", \"58\u00f75=11, 3\"],\n  [\"30\u00f79=3, 3\", \"15\u00f76=2, 3\"],\n  [\"51\u00f79=5, 6\", \"90\u00f75=18, 0\"],\n  [\"28\u00f73=9, 1\", \"64\u00f76=10, 4\"],\n  [\"49\u00f73=16, 1\", \"31\u00f77=4, 3\"],\n  [\"90\u00f79=10, 0\", \"12\u00f76=2, 0\"],\n  [\"60\u00f74=15, 0\", \"49\u00f77=7, 0\"],\n  [\"20\u00f73=6, 2\", \"44\u00f74=11, 0\"],\n  [\"35\u00f77=5, 0\", \"73\u00f78=9, 1\"],\n  [\"55\u00f73=18, 1\", \"49\u00f74=12, 1\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and the 25 division-problem answers in the table.\n# Each \"before\" text is unique within the document, so Find/Replace keyed\n# on exact text is unambiguous and safe regardless of order (no \"after\"\n# value collides with any other \"before\" value).\n#\n# Note: wdFindContinue / wdReplaceAll named constants are not predefined\n# in this environment, so the literal values (1 and 2 respectively) are\n# used directly with Find.Execute.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-03 Saturday\", \"2025-05-04 Sunday\"),\n    @(\"87\u00f72=43, 1\", \"51\u00f72=25, 1\"),\n    @(\"33\u00f79=3, 6\", \"87\u00f79=9, 6\"),\n    @(\"51\u00f75=10, 1\", \"16\u00f79=1, 7\"),\n    @(\"47\u00f72=23, 1\", \"36\u00f78=4, 4\"),\n    @(\"97\u00f72=48, 1\", \"25\u00f76=4, 1\"),\n    @(\"71\u00f72=35, 1\", \"10\u00f78=1, 2\"),\n    @(\"71\u00f75=14, 1\", \"18\u00f75=3, 3\"),\n    @(\"55\u00f76=9, 1\", \"96\u00f75=19, 1\"),\n    @(\"99\u00f75=19, 4\", \"38\u00f77=5, 3\"),\n    @(\"40\u00f78=5, 0\", \"13\u00f79=1, 4\"),\n    @(\"58\u00f72=29, 0\", \"23\u00f75=4, 3\"),\n    @(\"13\u00f78=1, 5\", \"82\u00f73=27, 1\"),\n    @(\"53\u00f76=8, 5\", \"17\u00f75=3, 2\"),\n    @(\"35\u00f72=17, 1\", \"31\u00f73=10, 1\"),\n    @(\"79\u00f77=11, 2\", \"52\u00f75=10, 2\"),\n    @(\"36\u00f75=7, 1\", \"58\u00f75=11, 3\"),\n    @(\"30\u00f79=3, 3\", \"15\u00f76=2, 3\"),\n    @(\"51\u00f79=5, 6\", \"90\u00f75=18, 0\"),\n    @(\"28\u00f73=9, 1\", \"64\u00f76=10, 4\"),\n    @(\"49\u00f73=16, 1\", \"31\u00f77=4, 3\"),\n    @(\"90\u00f79=10, 0\", \"12\u00f76=2, 0\"),\n    @(\"60\u00f74=15, 0\", \"49\u00f77=7, 0\"),\n    @(\"20\u00f73=6, 2\", \"44\u00f74=11, 0\"),\n    @(\"35\u00f77=5, 0\", \"73\u00f78=9, 1\"),\n    @(\"55\u00f73=18, 1\", \"49\u00f74=12, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    # 1 = wdFindContinue, 2 = wdReplaceAll\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
